$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column F: "care_home_deaths" -------------------------------------
# Header first (matches the order new shared strings appear in the target
# file: "DNK1_nch" was typed into A19 before the header went in).
$ws.Range("A19").Value = "DNK1_nch"

$ws.Range("F1").Value = "care_home_deaths"
# Copy the header formatting (bold, left aligned) from an existing header
# cell so we reuse the existing style (s="1") instead of minting a new one.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Existing rows (2-18) all had care home deaths included -> "yes"
$ws.Range("F2").Value = "yes"
$ws.Range("F2:F18").Value = "yes"

# New rows (19-24) are the "no care home" (noCH) variants -> "no"
$ws.Range("F19").Value = "no"

# --- New study_id values in column A for rows 20-24 ------------------------
$ws.Range("A20").Value = "ESP1-2_nch"
$ws.Range("A21").Value = "GBR3_nch"
$ws.Range("A22").Value = "CHE1_nch"
$ws.Range("A23").Value = "CHE2_nch"
$ws.Range("A24").Value = "NYC_NY_1_nch"

# --- New relpath values in column C for rows 19-24 --------------------------
$ws.Range("C19").Value = "data/derived/DNK1/DNK1_agebands_noCH.RDS"
$ws.Range("C20").Value = "data/derived/ESP1-2/ESP1-2_agebands_noCH.RDS"
$ws.Range("C21").Value = "data/derived/GBR3/GBR3_agebands_noCH.RDS"
$ws.Range("C22").Value = "data/derived/CHE1/CHE1_agebands_noCH.RDS"
$ws.Range("C23").Value = "data/derived/CHE2/CHE2_agebands_noCH.RDS"
$ws.Range("C24").Value = "data/derived/USA/NYC_NY_1_agebands_noCH.RDS"

# --- Fill in the remaining breakdown / serology_type / death_type columns --
$ws.Range("B19").Value = "ageband"
$ws.Range("D19").Value = "marginal"
$ws.Range("E19").Value = "aggregate"

$ws.Range("B20").Value = "ageband"
$ws.Range("D20").Value = "marginal"
$ws.Range("E20").Value = "aggregate"

$ws.Range("B21").Value = "ageband"
$ws.Range("D21").Value = "marginal"
$ws.Range("E21").Value = "aggregate"

$ws.Range("B22").Value = "ageband"
$ws.Range("D22").Value = "stratified"
$ws.Range("E22").Value = "aggregate"

$ws.Range("B23").Value = "ageband"
$ws.Range("D23").Value = "stratified"
$ws.Range("E23").Value = "aggregate"

$ws.Range("B24").Value = "ageband"
$ws.Range("D24").Value = "marginal"
$ws.Range("E24").Value = "aggregate"

# --- Selection / scroll position as left by the editor ---------------------
$ws.Range("C25").Select()
